# Edit script for LOT2015.xlsx
# Fixes misaligned content in the "Engenharia Bioquimica II" discipline sheet:
# inserts a dedicated row for the responsible-teacher entry, and fills in the
# objectives / short-summary / full-program / method / criteria / recovery /
# bibliography texts that were missing or misplaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long text values -------------------------------------------------
$NEW_OBJECTIVES = "Capacitar os alunos para a aplicação de conhecimentos da engenharia química na solução de problemas que se apresentam na implantação e otimização de processos biotecnológicos, com ênfase em: 1) agitação e aeração em processos fermentativos; 2) ampliação de escala e; 3) recuperação de produtos biotecnológicos."
$NEW_RESUMO = "A importância da transferência de oxigênio; sistemas de transferência de oxigênio; transferência de oxigênio e respiração microbiana; transferência de oxigênio em sistemas agitados e aerados, variação de escala, purificação de produtos biotecnológicos."
$NEW_PROGRAMA = "1. A importância da transferência de oxigênio.2. Sistemas de transferência de oxigênio.3. Transferência de oxigênio e respiração microbiana: análise conjunta da transferência e do consumo de oxigênio, determinação de kLa e de QO2 durante o processo fermentativo.4. Transferência de oxigênio em sistemas agitados e aerados: agitação de líquidos Newtonianos, agitação de líquidos Newtonianos submetidos à aeração, agitação de líquidos não Newtonianos, transferência de oxigênio.5. Variação de escala: critérios para ampliação de escala, comparações entre critérios para aampliação de escala, redução de escala.6. Introdução à purificação de produtos biotecnológicos: clarificação, rompimento de células, técnicas de separação de bioprodutos."
$METODO_TEXT = "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2."
$CRITERIO_TEXT = "A nota final (NF) será calculada como: N_F=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$NORMA_TEXT = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR).`nA média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$NEW_BIBLIO = "1. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial - Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001.`n2. SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. Biotecnologia Industrial - Engenharia Bioquímica (Vol 2), São Paulo: Edgard Blucher Ltda, 2001.`n3. DORAN P.M. Bioprocess Engineering Principles, 1st edition, San Diego: Academic Press, 1995.`n4. BAILEY, J.E., OLLIS D.F. Biochemical Engineering Fundamentals. 2nd edition, New York: McGraw Hill, 1986."
$TEACHER = "101761 - Arnaldo Márcio Ramalho Prata"

# --- Row 10: Objetivos / Objectives content ---------------------------
$ws.Range("B10").Value = $NEW_OBJECTIVES
$ws.Range("C10").Value = $NEW_OBJECTIVES

# --- Insert a new row at 13 for the "Docentes responsaveis" entry -----
# (previously the teacher name lived on the "Objetivos" row; it now gets
# its own row right below "Docentes responsaveis:")
$ws.Rows.Item(13).Insert()

$ws.Range("B13").Value = $TEACHER
$ws.Range("C13").Value = $TEACHER

# --- Row 14: Programa resumido -----------------------------------------
$ws.Range("B14").Value = $NEW_RESUMO
$ws.Range("C14").Value = $NEW_RESUMO

# --- Row 16: Programa (full PT syllabus) -------------------------------
$ws.Range("B16").Value = $NEW_PROGRAMA
$ws.Range("C16").Value = $NEW_PROGRAMA

# --- Row 19: Metodo ------------------------------------------------------
$ws.Range("B19").Value = $METODO_TEXT
$ws.Range("C19").Value = $METODO_TEXT

# --- Row 20: Criterio ------------------------------------------------------
$ws.Range("B20").Value = $CRITERIO_TEXT
$ws.Range("C20").Value = $CRITERIO_TEXT

# --- Row 21: Norma de recuperacao ------------------------------------------
$ws.Range("B21").Value = $NORMA_TEXT
$ws.Range("C21").Value = $NORMA_TEXT

# --- Row 22: Bibliografia --------------------------------------------------
$ws.Range("B22").Value = $NEW_BIBLIO
$ws.Range("C22").Value = $NEW_BIBLIO

# --- Column layout cleanup: column A should only cover col 1, not 1:2 -----
$ws.Columns.Item(2).Hidden = $false
